# Update crypto price/volume figures per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.673.56'
$ws.Range("E2").Value = '  -1.72%  '
$ws.Range("D3").Value = '2.905.17'
$ws.Range("E3").Value = '  -2.17%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'528.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.45%  '
$ws.Range("D6").Value = "'144.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.68%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'0.557"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.27%  '
$ws.Range("D9").Value = '2.913.25'
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("E10").Value = '  -3.51%  '
$ws.Range("D11").Value = "'6.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = '3.411.23'
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("E14").Value = '  +2.53%  '
$ws.Range("D15").Value = '60.648.06'
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").Value = "'22.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.67%  '
$ws.Range("D17").Value = '2.907.54'
$ws.Range("E17").Value = '  -2.22%  '
$ws.Range("E18").Value = '  -3.69%  '
$ws.Range("D19").Value = "'5.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.20%  '
$ws.Range("D20").Value = "'11.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.86%  '
$ws.Range("D21").Value = "'362.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.53%  '
$ws.Range("D22").Value = "'6.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = "'5.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("D25").Value = "'64.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.41%  '
$ws.Range("E26").Value = '  -2.83%  '
$ws.Range("E27").Value = '  -2.90%  '
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").Value = "'7.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.16%  '
$ws.Range("D30").Value = '0.0₃0863'
$ws.Range("E30").Value = '  -7.49%  '
$ws.Range("E32").Value = '  -2.31%  '
$ws.Range("D33").Value = "'19.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.46%  '
$ws.Range("D34").Value = "'152.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.94%  '
$ws.Range("D35").Value = "'4.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.50%  '
$ws.Range("D36").Value = "'5.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.72%  '
$ws.Range("E37").Value = '  -4.59%  '
$ws.Range("E38").Value = '  -5.52%  '
$ws.Range("D39").Value = "'37.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("E40").Value = '  -4.48%  '
$ws.Range("D41").Value = "'3.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.61%  '
$ws.Range("D42").Value = '2.295.36'
$ws.Range("E42").Value = '  -4.80%  '
$ws.Range("D43").Value = "'0.649"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.58%  '
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("D45").Value = "'20.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.06%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("E47").Value = '  +0.96%  '
$ws.Range("E48").Value = '  -2.88%  '
$ws.Range("E49").Value = '  -1.44%  '
$ws.Range("E50").Value = '  -2.14%  '
$ws.Range("D51").Value = "'252.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.21%  '
